# Weekly update: a new Jengibre price record for Mercado Mayorista Lo
# Valledor de Santiago (week of 2023-08-14) is inserted as a new data row
# right after the header/existing rows 2-41, pushing all the subsequent
# historical rows down by one. The sheet's used range grows from
# A1:R143 to A1:R144.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 42 - this shifts rows 42..143 down to
# 43..144 (and the whole sheet dimension grows to A1:R144) while keeping
# all of their existing values/styles intact, exactly like Excel's own
# "Insert Sheet Rows" command.
$ws.Rows.Item(42).Insert()

# Populate the freshly inserted row 42 with the new weekly observation.
$ws.Cells.Item(42, 1).Value = 6
$ws.Cells.Item(42, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(42, 3).Value = "Metropolitana"
$ws.Cells.Item(42, 4).Value = 45152
$ws.Cells.Item(42, 5).Value = 13
$ws.Cells.Item(42, 6).Value = 100114007
$ws.Cells.Item(42, 7).Value = "Jengibre"
$ws.Cells.Item(42, 8).Value = "Sin especificar"
$ws.Cells.Item(42, 9).Value = "Primera"
$ws.Cells.Item(42, 10).Value = 220
$ws.Cells.Item(42, 11).Value = 15000
$ws.Cells.Item(42, 12).Value = 17000
$ws.Cells.Item(42, 13).Value = 15909
$ws.Cells.Item(42, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(42, 15).Value = "Perú"
$ws.Cells.Item(42, 16).Value = 1224
$ws.Cells.Item(42, 17).Value = 13
$ws.Cells.Item(42, 18).Value = "Hortaliza"
